$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: the review flag for this entry changes from "yes" to "no"
$ws.Range("G6").Value = "no"

# New row 10: additional submission for com.singleton.strechy / taxi
# Clone row 9's formatting first (keeps styles.xml untouched, matching the
# existing per-column styles: A=1, B=0, C=2, D=2, E=0, F=0, G=0) then
# overwrite the values for the new record (no F10 value -> no review text).
$ws.Range("A9:G9").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)

$ws.Range("A10").Value = "com.singleton.strechy"
$ws.Range("B10").Value = "taxi"
$ws.Range("C10").Value = "cohenn167@gmail.com"
$ws.Range("D10").Value = "stavsade45@gmail.com"
$ws.Range("E10").Value = "27/5/2019 15:59"
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = "confirm"

# Update the sheet view: scrolled right one column, selection on G7
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("G7").Select()
